# Auto-generated script: apply scheduled-runner market-price updates to Raiden_Profits workbook
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 509.57144
$ws.Range("I33").Value = 509.57144
$ws.Range("K33").Value = 509.57144
$ws.Range("M33").Value = -280.57144
$ws.Range("H34").Value = 28077
$ws.Range("I34").Value = 7499
$ws.Range("J34").Value = 48655
$ws.Range("K34").Value = 7499
$ws.Range("L34").Value = 48655
$ws.Range("M34").Value = -7296
$ws.Range("N34").Value = -49061
$ws.Range("H36").Value = 28077
$ws.Range("I36").Value = 7499
$ws.Range("J36").Value = 48655
$ws.Range("K36").Value = 7499
$ws.Range("L36").Value = 48655
$ws.Range("M36").Value = -6784
$ws.Range("N36").Value = -50085
$ws.Range("H43").Value = 14906.833
$ws.Range("J43").Value = 36471.5
$ws.Range("L43").Value = 36471.5
$ws.Range("N43").Value = -36609.5
$ws.Range("H64").Value = 3798.3333
$ws.Range("J64").Value = 3798
$ws.Range("L64").Value = 3798
$ws.Range("N64").Value = -4294
$ws.Range("H67").Value = 3798.3333
$ws.Range("J67").Value = 3798
$ws.Range("L67").Value = 3798
$ws.Range("N67").Value = -5514
$ws.Range("H86").Value = 1941.2858
$ws.Range("I86").Value = 1549.5
$ws.Range("K86").Value = 1549.5
$ws.Range("M86").Value = -426.5
$ws.Range("H89").Value = 1941.2858
$ws.Range("I89").Value = 1549.5
$ws.Range("K89").Value = 7747.5
$ws.Range("M89").Value = -2131.5
$ws.Range("H92").Value = 588.7692
$ws.Range("I92").Value = 616.63635
$ws.Range("J92").Value = 435.5
$ws.Range("K92").Value = 616.63635
$ws.Range("L92").Value = 435.5
$ws.Range("M92").Value = 631.36365
$ws.Range("N92").Value = -2931.5
$ws.Range("H94").Value = 2316.1667
$ws.Range("I94").Value = 2316.1667
$ws.Range("K94").Value = 2316.1667
$ws.Range("M94").Value = -1865.1667
$ws.Range("H98").Value = 1293.25
$ws.Range("I98").Value = 1434.75
$ws.Range("J98").Value = 1010.25
$ws.Range("K98").Value = 1434.75
$ws.Range("L98").Value = 1010.25
$ws.Range("M98").Value = 63.25
$ws.Range("N98").Value = -4006.25
$ws.Range("H100").Value = 4941.3335
$ws.Range("I100").Value = 4895
$ws.Range("K100").Value = 4895
$ws.Range("M100").Value = -4354
$ws.Range("H122").Value = 1293.25
$ws.Range("I122").Value = 1434.75
$ws.Range("J122").Value = 1010.25
$ws.Range("K122").Value = 4304.25
$ws.Range("L122").Value = 3030.75
$ws.Range("M122").Value = -1854.25
$ws.Range("N122").Value = -7930.75
$ws.Range("H137").Value = 3272.647
$ws.Range("I137").Value = 2710.2666
$ws.Range("J137").Value = 3716.6316
$ws.Range("K137").Value = 8130.7998
$ws.Range("L137").Value = 11149.8948
$ws.Range("M137").Value = -5580.7998
$ws.Range("N137").Value = -16249.8948
$ws.Range("H138").Value = 2954.8135
$ws.Range("I138").Value = 2246.3635
$ws.Range("K138").Value = 6739.0905
$ws.Range("M138").Value = -1599.0905

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H46").Value = 18016.5
$ws.Range("I46").Value = 13849.75
$ws.Range("J46").Value = 26350
$ws.Range("K46").Value = 13849.75
$ws.Range("L46").Value = 26350
$ws.Range("M46").Value = -13530.75
$ws.Range("N46").Value = -26988
$ws.Range("H134").Value = 67500
$ws.Range("J134").Value = 67500
$ws.Range("L134").Value = 67500
$ws.Range("N134").Value = -77640

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2493.3333
$ws.Range("I86").Value = 2493.3333
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 2493.3333
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -1370.3333
$ws.Range("N86").ClearContents()
$ws.Range("H89").Value = 2493.3333
$ws.Range("I89").Value = 2493.3333
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 12466.6665
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = -6850.666499999999
$ws.Range("N89").ClearContents()

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 83479.586
$ws.Range("I7").Value = 200151.8
$ws.Range("J7").Value = 142.28572
$ws.Range("K7").Value = 200151.8
$ws.Range("L7").Value = 142.28572
$ws.Range("M7").Value = -200038.8
$ws.Range("N7").Value = -368.28572
$ws.Range("H99").Value = 6499.25
$ws.Range("I99").Value = 6499.25
$ws.Range("K99").Value = 6499.25
$ws.Range("M99").Value = -5001.25
$ws.Range("H122").Value = 2190.9167
$ws.Range("I122").Value = 2434.6843
$ws.Range("K122").Value = 7304.0529
$ws.Range("M122").Value = -4854.0529
$ws.Range("H126").Value = 6499.25
$ws.Range("I126").Value = 6499.25
$ws.Range("K126").Value = 19497.75
$ws.Range("M126").Value = -17027.75

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 2965434.8
$ws.Range("I4").Value = 5189429.5
$ws.Range("K4").Value = 15568288.5
$ws.Range("M4").Value = -15568176.5
$ws.Range("H64").Value = 3302
$ws.Range("I64").Value = 4005.5
$ws.Range("J64").Value = 1895
$ws.Range("K64").Value = 12016.5
$ws.Range("L64").Value = 5685
$ws.Range("M64").Value = -11746.5
$ws.Range("N64").Value = -6225
$ws.Range("H67").Value = 3302
$ws.Range("I67").Value = 4005.5
$ws.Range("J67").Value = 1895
$ws.Range("K67").Value = 12016.5
$ws.Range("L67").Value = 5685
$ws.Range("M67").Value = -11080.5
$ws.Range("N67").Value = -7557

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H40").Value = 30018
$ws.Range("J40").Value = 30018
$ws.Range("L40").Value = 30018
$ws.Range("N40").Value = -30320
$ws.Range("H43").Value = 4236
$ws.Range("I43").Value = 2981.5557
$ws.Range("J43").Value = 7999.3335
$ws.Range("K43").Value = 2981.5557
$ws.Range("L43").Value = 7999.3335
$ws.Range("M43").Value = -2830.5557
$ws.Range("N43").Value = -8301.333500000001
$ws.Range("H132").Value = 3563.6553
$ws.Range("I132").Value = 3563.6553
$ws.Range("K132").Value = 10690.9659
$ws.Range("M132").Value = -8160.965899999999

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1077.125
$ws.Range("I16").Value = 1077.125
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 1077.125
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -907.125
$ws.Range("N16").ClearContents()
$ws.Range("H22").Value = 657.2
$ws.Range("I22").Value = 684
$ws.Range("J22").Value = 550
$ws.Range("K22").Value = 684
$ws.Range("L22").Value = 550
$ws.Range("M22").Value = -389
$ws.Range("N22").Value = -1140
$ws.Range("H27").Value = 657.2
$ws.Range("I27").Value = 684
$ws.Range("J27").Value = 550
$ws.Range("K27").Value = 684
$ws.Range("L27").Value = 550
$ws.Range("M27").Value = -577
$ws.Range("N27").Value = -764
$ws.Range("H40").Value = 4181.25
$ws.Range("I40").Value = 4106.8184
$ws.Range("J40").Value = 5000
$ws.Range("K40").Value = 4106.8184
$ws.Range("L40").Value = 5000
$ws.Range("M40").Value = -3970.8184
$ws.Range("N40").Value = -5272
$ws.Range("H55").Value = 434
$ws.Range("I55").Value = 554.5454999999999
$ws.Range("J55").Value = 168.8
$ws.Range("K55").Value = 554.5454999999999
$ws.Range("L55").Value = 168.8
$ws.Range("M55").Value = -381.5454999999999
$ws.Range("N55").Value = -514.8
$ws.Range("H132").Value = 3224.4167
$ws.Range("I132").Value = 3307.3333
$ws.Range("K132").Value = 9921.999899999999
$ws.Range("M132").Value = -7391.999899999999

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H43").Value = 33332.332
$ws.Range("J43").Value = 19999
$ws.Range("L43").Value = 19999
$ws.Range("N43").Value = -20297
$ws.Range("H81").Value = 6163
$ws.Range("I81").Value = 4250
$ws.Range("K81").Value = 8500
$ws.Range("M81").Value = -7439
$ws.Range("H84").Value = 6163
$ws.Range("I84").Value = 4250
$ws.Range("K84").Value = 42500
$ws.Range("M84").Value = -37196
$ws.Range("H107").Value = 1031.75
$ws.Range("I107").Value = 374.77777
$ws.Range("J107").Value = 3002.6667
$ws.Range("K107").Value = 1124.33331
$ws.Range("L107").Value = 9008.000100000001
$ws.Range("M107").Value = 795.66669
$ws.Range("N107").Value = -12848.0001
$ws.Range("H113").Value = 761.125
$ws.Range("I113").Value = 761.125
$ws.Range("K113").Value = 2283.375
$ws.Range("M113").Value = -113.375
